# Updates the raw (non-formula) measurement inputs in columns BZ and CA
# for rows 17-44 of the "Measurements" sheet. All other changed cells in
# the original edit are formulas that depend on BZ/CA (directly or
# transitively) and will be recalculated automatically by Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Measurements")

# Scroll the window one column to the left (view previously started at BS5,
# now starts at BR5) without disturbing the current selection.
$excel.ActiveWindow.ScrollRow = 5
$excel.ActiveWindow.ScrollColumn = 70

# row -> (BZ, CA)
$updates = @{
    17 = @(451.14073305264702, 30.500284298828198)
    18 = @(451.03528844981503, 30.248306402218901)
    19 = @(447.42509295055402, 30.177805701306401)
    20 = @(441.66917251918397, 30.0942606097991)
    21 = @(440.02250817068102, 29.589694195267299)
    22 = @(434.38531765659201, 29.195317156250301)
    23 = @(423.51355177771597, 29.292151475042299)
    24 = @(417.478031635435, 37.863349149924801)
    25 = @(412.01912364321902, 36.930310328279297)
    26 = @(416.40573704738102, 36.809646026462197)
    27 = @(416.42646744760202, 37.166961182517198)
    28 = @(411.56563398078703, 37.221952432042301)
    29 = @(410.96820647045303, 37.335421002850801)
    30 = @(414.48973005282602, 36.813392980989697)
    31 = @(412.55946342898699, 38.583190778175798)
    32 = @(415.56523946965802, 38.352528287923001)
    33 = @(416.24947948158098, 37.4632658666658)
    34 = @(414.55828556110299, 37.265914065913201)
    35 = @(416.09284730552901, 37.417789833399702)
    36 = @(416.69894273247502, 36.520946379221002)
    37 = @(413.13338560641, 37.502743144990298)
    38 = @(415.27462013686102, 39.753063644963497)
    39 = @(418.05787780135898, 37.854383278473399)
    40 = @(411.12762566431599, 37.602762740977198)
    41 = @(418.05726359705102, 37.722479908920597)
    42 = @(417.54011305403299, 36.696877484428398)
    43 = @(418.90644549272201, 37.570939396349203)
    44 = @(416.25130174178997, 37.364548727306001)
}

foreach ($row in $updates.Keys) {
    $pair = $updates[$row]
    $ws.Cells.Item($row, 78).Value = $pair[0]   # column BZ
    $ws.Cells.Item($row, 79).Value = $pair[1]   # column CA
}
